$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '30.107.13'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '1.914.01'
$ws.Range('E3').Value = '  +0.35%  '
Set-TextCell 'D4' '1.001'
$ws.Range('E4').Value = '  +0.10%  '
Set-TextCell 'D5' '0.7973'
$ws.Range('E5').Value = '  +6.36%  '
Set-TextCell 'D6' '243.84'
$ws.Range('E6').Value = '  +0.28%  '
Set-TextCell 'D7' '1.000'
$ws.Range('E7').Value = '  +0.05%  '
Set-TextCell 'D8' '0.3191'
$ws.Range('E8').Value = '  +3.23%  '
Set-TextCell 'D9' '26.40'
$ws.Range('E9').Value = '  -0.32%  '
Set-TextCell 'D10' '0.06962'
$ws.Range('E10').Value = '  -0.21%  '
Set-TextCell 'D11' '0.08005'
$ws.Range('E11').Value = '  -1.03%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 'D12' '0.7528'
$ws.Range('E12').Value = '  -2.23%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.916.76'
$ws.Range('E13').Value = '  +0.28%  '
Set-TextCell 'D14' '5.247'
$ws.Range('E14').Value = '  -0.82%  '
Set-TextCell 'D15' '93.70'
$ws.Range('E15').Value = '  +1.73%  '
$ws.Range('D16').Value = '30.126.37'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('E17').Value = '  -0.87%  '
Set-TextCell 'D18' '5.968'
$ws.Range('E18').Value = '  -1.98%  '
Set-TextCell 'D19' '249.35'
$ws.Range('E19').Value = '  +3.82%  '
Set-TextCell 'D20' '0.000007831'
$ws.Range('E20').Value = '  +0.20%  '
Set-TextCell 'D21' '1.000'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 'D22' '1.000'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 'D23' '6.945'
$ws.Range('E23').Value = '  -2.75%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D24' '169.80'
$ws.Range('E24').Value = '  +1.63%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D25' '9.337'
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D26' '0.1402'
$ws.Range('E26').Value = '  +9.55%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 'D27' '18.99'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 'D28' '2.068'
$ws.Range('E28').Value = '  +0.78%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 'D29' '1.385'
$ws.Range('E29').Value = '  +2.52%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D30' '1.529'
$ws.Range('E30').Value = '  -0.59%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D31' '4.355'
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 'D32' '4.128'
$ws.Range('E32').Value = '  +1.23%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D33' '0.05530'
$ws.Range('E33').Value = '  +4.82%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 'D34' '1.267'
$ws.Range('E34').Value = '  -3.14%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D35' '0.7410'
$ws.Range('E35').Value = '  -0.74%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 'D36' '2.726'
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D37' '0.01933'
$ws.Range('E37').Value = '  -1.69%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 'D38' '2.799'
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D39' '6.218'
$ws.Range('E39').Value = '  -1.65%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 'D40' '0.4463'
$ws.Range('E40').Value = '  -0.56%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 'D41' '73.20'
$ws.Range('E41').Value = '  -1.34%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D42' '1.913'
$ws.Range('E42').Value = '  -3.07%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 'D43' '1.001'
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 'D44' '0.8339'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D45' '7.623'
$ws.Range('E45').Value = '  -1.30%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 'D46' '101.05'
$ws.Range('E46').Value = '  -0.84%  '
Set-TextCell 'D47' '9.877'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 'D48' '992.75'
$ws.Range('E48').Value = '  +7.50%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.069.05'
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell 'D50' '36.58'
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D51' '1.509'
$ws.Range('E51').Value = '  +1.22%  '
